$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.019.99"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.679.02"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'215.80"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "'21.40"
$ws.Range("E9").Value = "  +5.47%  "
$ws.Range("D10").Value = "'0.0624"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "1.915.66"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "1.695.08"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'66.53"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "27.022.66"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'8.15"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "'235.75"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").Value = "'146.43"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "'16.42"
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "1.539.31"
$ws.Range("E33").Value = "  +5.52%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").Value = "'1.72"
$ws.Range("E35").Value = "  +5.15%  "
$ws.Range("D36").Value = "'2.39"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "'0.920"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("E39").Value = "  +3.19%  "
$ws.Range("E40").Value = "  +6.44%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'67.96"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("D43").Value = "'5.60"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "1.821.57"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'90.36"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("E50").Value = "  +5.58%  "
$ws.Range("E51").Value = "  -0.01%  "

Write-Host "Applied crypto price/volume updates"
